$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.946.85'
$ws.Range("E2").Value = '  -4.54%  '
$ws.Range("D3").Value = '2.915.32'
$ws.Range("E3").Value = '  -7.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '475.13'
$ws.Range("E5").Value = '  -9.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.44'
$ws.Range("E6").Value = '  -3.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '2.913.87'
$ws.Range("E8").Value = '  -7.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.410'
$ws.Range("E9").Value = '  -8.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.72'
$ws.Range("E10").Value = '  -6.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0984'
$ws.Range("E11").Value = '  -11.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.338'
$ws.Range("E12").Value = '  -13.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.124'
$ws.Range("E13").Value = '  -2.87%  '
$ws.Range("D14").Value = '3.414.32'
$ws.Range("E14").Value = '  -7.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.66'
$ws.Range("E15").Value = '  -8.61%  '
$ws.Range("D16").Value = '54.910.11'
$ws.Range("E16").Value = '  -4.60%  '
$ws.Range("D17").Value = '2.911.36'
$ws.Range("E17").Value = '  -7.34%  '
$ws.Range("E18").Value = '  -11.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.46'
$ws.Range("E19").Value = '  -5.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.60'
$ws.Range("E20").Value = '  -11.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.17'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '305.02'
$ws.Range("E22").Value = '  -12.20%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.448'
$ws.Range("E24").Value = '  -12.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '59.59'
$ws.Range("E25").Value = '  -14.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.994'
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.153'
$ws.Range("E27").Value = '  -7.71%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = '0.0₃0817'
$ws.Range("E29").Value = '  -14.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.21'
$ws.Range("E30").Value = '  -9.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.33'
$ws.Range("E31").Value = '  -8.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.13'
$ws.Range("E32").Value = '  -7.05%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.65'
$ws.Range("E33").Value = '  -12.37%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.96'
$ws.Range("E34").Value = '  -12.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '144.71'
$ws.Range("E35").Value = '  -9.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.26'
$ws.Range("E36").Value = '  -13.89%  '
$ws.Range("E37").Value = '  -12.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.23'
$ws.Range("E38").Value = '  -12.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.37'
$ws.Range("E39").Value = '  -10.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0628'
$ws.Range("E40").Value = '  -9.70%  '
$ws.Range("D41").Value = '2.943.06'
$ws.Range("E41").Value = '  -7.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.63'
$ws.Range("E43").Value = '  -11.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.967'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.613'
$ws.Range("E45").Value = '  -11.43%  '
$ws.Range("E46").Value = '  -8.49%  '
$ws.Range("E47").Value = '  -12.59%  '
$ws.Range("D48").Value = '2.062.03'
$ws.Range("E48").Value = '  -8.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.41'
$ws.Range("E49").Value = '  -13.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0219'
$ws.Range("E50").Value = '  -6.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.10'
$ws.Range("E51").Value = '  -11.73%  '
